# Petalburg Gym trainers — add TRAINER_BERKE, TRAINER_JODY and
# TRAINER_NORMAN_1 party blocks right after the last existing trainer
# (mirrors the upstream "petalburg gym trainers, fix porymap" commit).
#
# NOTE: in the source workbook the old "END" marker cell (A94) is left
# pointing at a shared-string index that the newly-inserted strings push
# forward, so after this edit A94 reads "TRAINER_BERKE" and a new "END"
# sentinel row is appended at the bottom (A111) — that is the upstream
# porymap quirk this commit is fixing, and we reproduce it verbatim by
# simply appending the new blocks the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TRAINER_BERKE -----------------------------------------------------
$ws.Cells.Item(94, 1).Value = "TRAINER_BERKE"

$ws.Cells.Item(95, 1).Value = "species"
$ws.Cells.Item(95, 2).Value = "lvl"
$ws.Cells.Item(95, 3).Value = "iv"
$ws.Cells.Item(95, 4).Value = "heldItem"
$ws.Cells.Item(95, 5).Value = "moves"

$ws.Cells.Item(96, 1).Value = "Sentret"
$ws.Cells.Item(96, 2).Value = 8
$ws.Cells.Item(96, 5).Value = "Quick Attack"

$ws.Cells.Item(97, 1).Value = "Rattata"
$ws.Cells.Item(97, 2).Value = 10
$ws.Cells.Item(97, 4).Value = "Oran Berry"
$ws.Cells.Item(97, 5).Value = "Quick Attack"

# --- TRAINER_JODY --------------------------------------------------------
$ws.Cells.Item(99, 1).Value = "TRAINER_JODY"

$ws.Cells.Item(100, 1).Value = "species"
$ws.Cells.Item(100, 2).Value = "lvl"
$ws.Cells.Item(100, 3).Value = "iv"
$ws.Cells.Item(100, 4).Value = "heldItem"
$ws.Cells.Item(100, 5).Value = "moves"

$ws.Cells.Item(101, 1).Value = "Meowsy"
$ws.Cells.Item(101, 2).Value = 9
$ws.Cells.Item(101, 4).Value = "Oran Berry"
$ws.Cells.Item(101, 5).Value = "Pay Day"

$ws.Cells.Item(102, 1).Value = "Zigzagoon"
$ws.Cells.Item(102, 2).Value = 10
$ws.Cells.Item(102, 5).Value = "Headbutt"

# --- TRAINER_NORMAN_1 ------------------------------------------------------
$ws.Cells.Item(104, 1).Value = "TRAINER_NORMAN_1"

$ws.Cells.Item(105, 1).Value = "species"
$ws.Cells.Item(105, 2).Value = "lvl"
$ws.Cells.Item(105, 3).Value = "iv"
$ws.Cells.Item(105, 4).Value = "heldItem"
$ws.Cells.Item(105, 5).Value = "moves"

$ws.Cells.Item(106, 1).Value = "Galarian_Zigzagoon"
$ws.Cells.Item(106, 2).Value = 12
$ws.Cells.Item(106, 5).Value = "Snarl, Headbutt, Sand-Attack, Leer"

$ws.Cells.Item(107, 1).Value = "Vigoroth"
$ws.Cells.Item(107, 2).Value = 14
$ws.Cells.Item(107, 5).Value = "Scratch, Uproar, Fury Swipes, Focus Energy"

# --- new END sentinel ------------------------------------------------------
$ws.Cells.Item(111, 1).Value = "END"

# Move the selection/view to where the author ended up after typing the
# new block, as recorded in the sheetView.
$ws.Range("E108").Select()
